$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H92").Value = 842
$ws.Range("I92").Value = 974.2857
$ws.Range("J92").Value = 533.3333
$ws.Range("K92").Value = 974.2857
$ws.Range("L92").Value = 533.3333
$ws.Range("M92").Value = 273.7143
$ws.Range("N92").Value = -3029.3333
$ws.Range("H99").Value = 197.33333
$ws.Range("I99").Value = 196
$ws.Range("J99").Value = 200
$ws.Range("K99").Value = 588
$ws.Range("L99").Value = 600
$ws.Range("M99").Value = 910
$ws.Range("N99").Value = -3596
$ws.Range("H125").Value = 457.5625
$ws.Range("I125").Value = 372.66666
$ws.Range("K125").Value = 3353.99994
$ws.Range("M125").Value = -893.9999399999997
$ws.Range("H129").Value = 304189.66
$ws.Range("J129").Value = 386006.97
$ws.Range("L129").Value = 1158020.91
$ws.Range("N129").Value = -1168020.91
$ws.Range("H138").Value = 3487.389
$ws.Range("I138").Value = 7300
$ws.Range("J138").Value = 3140.7878
$ws.Range("K138").Value = 21900
$ws.Range("L138").Value = 9422.3634
$ws.Range("M138").Value = -16760
$ws.Range("N138").Value = -19702.3634

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 3220
$ws.Range("I28").Value = 3220
$ws.Range("K28").Value = 3220
$ws.Range("M28").Value = -3028
$ws.Range("H32").Value = 5309.708
$ws.Range("I32").Value = 4007.939
$ws.Range("J32").Value = 20559
$ws.Range("K32").Value = 4007.939
$ws.Range("L32").Value = 20559
$ws.Range("M32").Value = -3720.939
$ws.Range("N32").Value = -21133
$ws.Range("I45").Value = 3313.6667
$ws.Range("J45").Value = 2435.2727
$ws.Range("K45").Value = 3313.6667
$ws.Range("L45").Value = 2435.2727
$ws.Range("M45").Value = -2936.6667
$ws.Range("N45").Value = -3189.2727
$ws.Range("H99").Value = 3220
$ws.Range("I99").Value = 3220
$ws.Range("K99").Value = 3220
$ws.Range("M99").Value = -225

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 862.25
$ws.Range("I22").Value = 871.1429000000001
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 871.1429000000001
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = -698.1429000000001
$ws.Range("N22").Value = -1146
$ws.Range("H86").Value = 1915
$ws.Range("I86").Value = 1763.421
$ws.Range("J86").Value = 2275
$ws.Range("K86").Value = 1763.421
$ws.Range("L86").Value = 2275
$ws.Range("M86").Value = -640.421
$ws.Range("N86").Value = -4521
$ws.Range("H89").Value = 1915
$ws.Range("I89").Value = 1763.421
$ws.Range("J89").Value = 2275
$ws.Range("K89").Value = 8817.105
$ws.Range("L89").Value = 11375
$ws.Range("M89").Value = -3201.105
$ws.Range("N89").Value = -22607
$ws.Range("H94").Value = 1131.1428
$ws.Range("I94").Value = 956.5714
$ws.Range("J94").Value = 1305.7142
$ws.Range("K94").Value = 956.5714
$ws.Range("L94").Value = 1305.7142
$ws.Range("M94").Value = -505.5714
$ws.Range("N94").Value = -2207.7142
$ws.Range("H134").Value = 4896.48
$ws.Range("I134").Value = 5037.5835
$ws.Range("J134").Value = 1510
$ws.Range("K134").Value = 15112.7505
$ws.Range("L134").Value = 4530
$ws.Range("M134").Value = -12577.7505
$ws.Range("N134").Value = -9600

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4501.231
$ws.Range("I132").Value = 3172.9092
$ws.Range("K132").Value = 9518.7276
$ws.Range("M132").Value = -6988.7276
$ws.Range("H141").Value = 25274.447
$ws.Range("J141").Value = 25274.447
$ws.Range("L141").Value = 25274.447
$ws.Range("N141").Value = -35634.447

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H68").Value = 1166.3334
$ws.Range("J68").Value = 1250
$ws.Range("L68").Value = 3750
$ws.Range("N68").Value = -5372
$ws.Range("H71").Value = 1166.3334
$ws.Range("J71").Value = 1250
$ws.Range("L71").Value = 11250
$ws.Range("N71").Value = -19362
$ws.Range("H97").Value = 812.375
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 812.375
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 2437.125
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -3429.125
$ws.Range("H131").Value = 711.25
$ws.Range("J131").Value = 713.6869
$ws.Range("L131").Value = 2141.0607
$ws.Range("N131").Value = -12221.0607

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H80").Value = 3244.2693
$ws.Range("I80").Value = 2843
$ws.Range("K80").Value = 2843
$ws.Range("M80").Value = -1845
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H83").Value = 3244.2693
$ws.Range("I83").Value = 2843
$ws.Range("K83").Value = 14215
$ws.Range("M83").Value = -9223
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H113").Value = 4421.5
$ws.Range("I113").Value = 5422.15
$ws.Range("J113").Value = 1919.875
$ws.Range("K113").Value = 5422.15
$ws.Range("L113").Value = 1919.875
$ws.Range("M113").Value = -3252.15
$ws.Range("N113").Value = -6259.875
$ws.Range("H141").Value = 84479.8
$ws.Range("J141").Value = 84479.8
$ws.Range("L141").Value = 84479.8
$ws.Range("N141").Value = -94839.8

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1340.2
$ws.Range("I136").Value = 1385.9286
$ws.Range("J136").Value = 700
$ws.Range("K136").Value = 4157.7858
$ws.Range("L136").Value = 2100
$ws.Range("M136").Value = -1607.7858
$ws.Range("N136").Value = -7200

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1011.3
$ws.Range("I132").Value = 820.9231
$ws.Range("J132").Value = 1364.8572
$ws.Range("K132").Value = 2462.7693
$ws.Range("L132").Value = 4094.5716
$ws.Range("M132").Value = 67.23070000000007
$ws.Range("N132").Value = -9154.571599999999
$ws.Range("H140").Value = 46900
$ws.Range("J140").Value = 46900
$ws.Range("L140").Value = 46900
$ws.Range("N140").Value = -57260
$ws.Range("H141").Value = 80357.5
$ws.Range("J141").Value = 80357.5
$ws.Range("L141").Value = 80357.5
$ws.Range("N141").Value = -90717.5
